$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired state of A1:C19 after the row reorder described in the diff.
$data = @(
    @("Oyuncu Adı", "Pozisyon", "Takım"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Rui Hachimura", "SF,PF", "Los Angeles Lakers"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Kyle Kuzma", "PF", "Washington Wizards"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($rowNum, 1).Value = $row[0]
    $ws.Cells.Item($rowNum, 2).Value = $row[1]
    $ws.Cells.Item($rowNum, 3).Value = $row[2]
}
